$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Change 1: merge the "(en particulier : ... publique" run with the
#           ") " run that follows the old "_GoBack" bookmark into a
#           single run, dropping that now-stale bookmark.
# ---------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("(en particulier")
$mergeStart = $r1.Start

$r2 = $d.Content
$r2.Find.Execute(")" + [char]32)
$mergeEnd = $r2.End

# Temporary bookmarks shield the neighbouring runs (which share the
# same rPr) from being swept into the same merge as the edited range.
$d.Bookmarks.Add("zzProtectBefore", $d.Range($mergeStart, $mergeStart))
$d.Bookmarks.Add("zzProtectAfter", $d.Range($mergeEnd, $mergeEnd))

$mergeRange = $d.Range($mergeStart, $mergeEnd)
$mergeRange.Text = "(en particulier : ceux liés aux mesures de santé publique) "

$d.Bookmarks("zzProtectBefore").Delete()
$d.Bookmarks("zzProtectAfter").Delete()

# ---------------------------------------------------------------
# Change 2: italicise "CovidÉcoute" and relocate the "_GoBack"
#           bookmark onto it (Word auto-tracks the last edit point;
#           adding a bookmark named "_GoBack" replaces the old one).
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("CovidÉcoute", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Italic = 1
$d.Bookmarks.Add("_GoBack", $rng)
